# Weekly data refresh: insert a new observation row above row 66
# (pushing the existing rows 66-85 down to 67-86) and populate it with
# the latest "Camote" price record for Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 66; existing rows 66-85 shift to 67-86
# and the sheet's used range / dimension grows to A1:R86 automatically.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value = 9
$ws.Range("B66").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C66").Value = "Metropolitana"
$ws.Range("D66").Value = 44711
$ws.Range("E66").Value = 13
$ws.Range("F66").Value = 100114002
$ws.Range("G66").Value = "Camote"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 970
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14500
$ws.Range("N66").Value = "$/malla 18 kilos"
$ws.Range("O66").Value = "Perú"
$ws.Range("P66").Value = 806
$ws.Range("Q66").Value = 18
$ws.Range("R66").Value = "Hortaliza"
